# Publish terminology IG 2.0.0
#
# 1) Split the old "Concepts" sheet into two sheets:
#      - "Properties" (same position/sheetId as the old "Concepts")
#      - "Concepts"   (new sheet, holding what used to be in "Concepts")
# 2) Refresh the IG metadata (Version + Date) on the "Metadata" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Split Concepts -> Properties (old data) + Concepts (new sheet) ----
$concepts = $wb.Worksheets.Item("Concepts")

# Duplicate the sheet (keeps formatting/styles identical - no new style
# entries are introduced) and place the duplicate right after it.
$concepts.Copy([System.Reflection.Missing]::Value, $concepts)

# The original keeps its tab position/sheetId, but becomes "Properties".
$concepts.Name = "Properties"

# The duplicate becomes the "new" Concepts tab (unchanged content).
$copy = $wb.Worksheets.Item("Concepts (2)")
$copy.Name = "Concepts"

# --- 2. Overwrite "Properties" with the FHIR CodeSystem properties table --
$props = $wb.Worksheets.Item("Properties")

$props.Range("A1").Value = "Code"
$props.Range("B1").Value = "Uri"
$props.Range("C1").Value = "Description"
$props.Range("D1").Value = "Type"

$props.Range("A2").Value = "status"
$props.Range("B2").Value = "http://hl7.org/fhir/concept-properties#status"
$props.Range("C2").Value = "A property that indicates the status of the concept. One of active, experimental, deprecated, or retired."
$props.Range("D2").Value = "code"

$props.Range("A3").Value = "effectiveDate"
$props.Range("B3").Value = "http://hl7.org/fhir/concept-properties#effectiveDate"
$props.Range("C3").Value = "The date at which the concept status was last changed."
$props.Range("D3").Value = "dateTime"

# The old "Concepts" data occupied 4 rows; "Properties" only needs 3, so
# drop the now-unused trailing row entirely (keeps dimension = A1:D3).
$props.Rows.Item(4).Delete()

# --- 3. Bump the published IG version / date on the Metadata sheet -------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "1.2.1"

# Pre-format as Text so Excel doesn't reinterpret the ISO date string as a
# serial date number - the sheet stores it as a literal string.
$meta.Range("B8").NumberFormat = "@"
$meta.Range("B8").Value = "2025-09-22"
